$wb = $excel.ActiveWorkbook

# --- Sheet 2: "R-Wert und 7-Tage-Inzidenz" ---
# C6/C8 hold digit-only text ("184"/"334") that must stay text (not become
# numbers). Force text via NumberFormat, set the value, then ClearFormats so
# no stray style index is left behind on the cell (matches the original,
# unstyled shared-string cell).
$wsR = $wb.Worksheets.Item("R-Wert und 7-Tage-Inzidenz")
$wsR.Range("C6").NumberFormat = "@"
$wsR.Range("C6").Value = "184"
$wsR.Range("C6").ClearFormats()
$wsR.Range("D6").Value = "11,5 %"
$wsR.Range("C8").NumberFormat = "@"
$wsR.Range("C8").Value = "334"
$wsR.Range("C8").ClearFormats()
$wsR.Range("D8").Value = " 9,2 %"

# --- Sheet 5: "Todesfälle und Fallsterblichkei" ---
$wsT = $wb.Worksheets.Item("Todesfälle und Fallsterblichkei")
$wsT.Range("B3").Value = "912 ( 4,7%)"
$wsT.Range("D3").Value = " 19,6%"
$wsT.Range("B5").Value = "3623 ( 2,8%)"
$wsT.Range("D5").Value = " 23,1%"

# --- Sheet 7: "Regionale Daten" ---
$wsD = $wb.Worksheets.Item("Regionale Daten")
$wsD.Range("D2").Value = 184
$wsD.Range("D4").Value = 174
$wsD.Range("D8").Value = 174
$wsD.Range("E8").Value = 25
$wsD.Range("D9").Value = 183
$wsD.Range("D12").Value = 148
$wsD.Range("D14").Value = 205
